# Weekly fruit/vegetable data update:
# Insert a new data row at row 45 (pushing existing rows 45-65 down to 46-66)
# and populate it with the new week's record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 45; this shifts rows 45:65 down to 46:66,
# carrying their formatting (including the date style in column D) with them.
$ws.Rows.Item(45).Insert()

# Populate the newly inserted row 45 with the new record's data.
$ws.Cells.Item(45, 1).Value  = 1
$ws.Cells.Item(45, 2).Value  = 'Agrícola del Norte S.A. de Arica'
$ws.Cells.Item(45, 3).Value  = 'Arica y Parinacota'
$ws.Cells.Item(45, 4).Value  = 44726
$ws.Cells.Item(45, 5).Value  = 15
$ws.Cells.Item(45, 6).Value  = 100112009
$ws.Cells.Item(45, 7).Value  = 'Acelga'
$ws.Cells.Item(45, 8).Value  = 'Sin especificar'
$ws.Cells.Item(45, 9).Value  = 'Primera'
$ws.Cells.Item(45, 10).Value = 200
$ws.Cells.Item(45, 11).Value = 2500
$ws.Cells.Item(45, 12).Value = 2800
$ws.Cells.Item(45, 13).Value = 2650
$ws.Cells.Item(45, 14).Value = '$/atado 2,5 a 3 kilos'
$ws.Cells.Item(45, 15).Value = 'Región de Arica y Parinacota'
$ws.Cells.Item(45, 16).Value = 883
$ws.Cells.Item(45, 17).Value = 3
$ws.Cells.Item(45, 18).Value = 'Hortaliza'
